# Commit message: "updated time when you updated gear"
#
# Two groups of rows in the gear-tracking sheet (grouped by the Class in
# column D) got their row order reshuffled because a player's "gear
# updated" timestamp (column A) changed, which moves that player's whole
# record to a different slot within their class group. Columns M/N are an
# unrelated Class/Count summary table that happens to share the same row
# numbers and must stay untouched.
#
# Group 1 - rows 22-24 (class Lahn):
#   before: 22=Misha, 23=WickedKiss, 24=zozrog
#   after:  22=WickedKiss, 23=zozrog, 24=Misha (Misha's gear-update time
#           refreshed to a later timestamp on the same day; the rest of
#           Misha's stats are unchanged)
#
# Group 2 - rows 35-38 (class Ninja):
#   before: 35=Mrclutchcakes, 36=Born8, 37=Maalinko174, 38=dennymane
#   after:  35=dennymane, 36=Mrclutchcakes, 37=Born8, 38=Maalinko174
#           (a straight rotation of the 4 full records, dennymane's
#           existing record moves to the top of the group)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-GearRow($row, $b, $c, $d, $a, $e, $f, $g, $h, $i, $j, $k, $l) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $j
    $ws.Cells.Item($row, 11).Value = $k
    $ws.Cells.Item($row, 12).Value = $l
}

# ---- Group 1: rows 22-24 (Lahn) ----

Set-GearRow 22 "WickedKiss" "Wicked_Kiss" "Lahn" `
    46047.33054398148 354 354 430 `
    "Offense, Defense, Flex, Elephant, Cannoneer" "No" "Yes" "Yes" 784

Set-GearRow 23 "zozrog" "Zozrog" "Lahn" `
    46032.86646990741 361 361 437 `
    "Offense, Flex, Cannoneer, Scout, Cross comms, Flag Placer" "Yes" "Yes" "Yes" 798

Set-GearRow 24 "Misha" "Mishailia" "Lahn" `
    46047.74747685185 376 336 429 `
    "Offense, Flex" "No" "No" "Yes" 805

# ---- Group 2: rows 35-38 (Ninja) ----

Set-GearRow 35 "dennymane" "Dennymane" "Ninja" `
    46045.00385416667 387 389 458 `
    "Offense" "Yes" "No" "Yes" 847

Set-GearRow 36 "Mrclutchcakes" "MrClutchCakes" "Ninja" `
    46044.99524305556 395 395 465 `
    "Offense, Defense, Flex, Cannoneer, Scout, Flag Placer, Shai" "No" "No" "Yes" 860

Set-GearRow 37 "Born8" "Nightmvre" "Ninja" `
    46040.06761574074 390 390 352 `
    "Offense, Flex" "No" "No" "Yes" 742

Set-GearRow 38 "Maalinko174" "Sparten" "Ninja" `
    46026.35682870371 373 375 451 `
    "Defense, Flex" "No" "No" "Yes" 826
